# "figuras sombreadas y dataset"
# Row 1 currently holds the text labels "SHAC_216".."SHAC_221" (shared
# strings). Replace them with the bare numeric well codes (216, 217, 220,
# 221) used by the rest of the dataset, which also empties the shared
# string table. Then move the active selection to D2 (first data cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 216
$ws.Range("B1").Value = 217
$ws.Range("C1").Value = 220
$ws.Range("D1").Value = 221

$ws.Range("D2").Select()
